$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "61.616.52"
$ws.Range("E2").Value = "  +0.94%  "

$ws.Range("D3").Value = "3.391.88"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue "D5" "576.85"
$ws.Range("E5").Value = "  +0.89%  "

Set-TextValue "D6" "137.17"
$ws.Range("E6").Value = "  +1.09%  "

$ws.Range("D8").Value = "3.391.84"
$ws.Range("E8").Value = "  +0.61%  "

$ws.Range("E9").Value = "  -0.50%  "

Set-TextValue "D10" "7.51"
$ws.Range("E10").Value = "  -0.95%  "

Set-TextValue "D11" "0.126"
$ws.Range("E11").Value = "  +2.46%  "

$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").Value = "3.966.61"
$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("E14").Value = "  +1.36%  "

$ws.Range("E15").Value = "  +2.05%  "

$ws.Range("D16").Value = "3.390.44"
$ws.Range("E16").Value = "  +0.78%  "

Set-TextValue "D17" "25.78"
$ws.Range("E17").Value = "  +2.51%  "

$ws.Range("D18").Value = "61.695.99"
$ws.Range("E18").Value = "  +0.83%  "

Set-TextValue "D19" "14.20"
$ws.Range("E19").Value = "  +1.55%  "

Set-TextValue "D20" "9.50"
$ws.Range("E20").Value = "  +1.02%  "

Set-TextValue "D21" "5.83"
$ws.Range("E21").Value = "  +0.46%  "

Set-TextValue "D22" "379.16"
$ws.Range("E22").Value = "  +1.47%  "

Set-TextValue "D23" "0.560"
$ws.Range("E23").Value = "  -1.33%  "

$ws.Range("D24").Value = "3.524.72"
$ws.Range("E24").Value = "  +0.57%  "

$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D26" "0.0000125"
$ws.Range("E26").Value = "  +7.21%  "

$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D27" "71.18"
$ws.Range("E27").Value = "  +1.04%  "

Set-TextValue "D28" "1.73"
$ws.Range("E28").Value = "  +3.30%  "

$ws.Range("E29").Value = "  -1.40%  "

Set-TextValue "D30" "1.01"
$ws.Range("E30").Value = "  +1.33%  "

$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D31" "0.160"
$ws.Range("E31").Value = "  +3.82%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D32" "8.18"
$ws.Range("E32").Value = "  +0.90%  "

$ws.Range("E33").Value = "  +1.02%  "

$ws.Range("E34").Value = "  +0.04%  "

Set-TextValue "D35" "23.39"
$ws.Range("E35").Value = "  +0.31%  "

Set-TextValue "D36" "5.34"
$ws.Range("E36").Value = "  -3.59%  "

$ws.Range("E37").Value = "  +0.27%  "

Set-TextValue "D38" "6.83"
$ws.Range("E38").Value = "  -1.09%  "

Set-TextValue "D39" "164.11"
$ws.Range("E39").Value = "  +0.19%  "

Set-TextValue "D40" "0.0782"
$ws.Range("E40").Value = "  -0.57%  "

$ws.Range("E41").Value = "  +2.91%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D42" "0.781"
$ws.Range("E42").Value = "  +2.84%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "1.00"
$ws.Range("E43").Value = "  -0.01%  "

Set-TextValue "D44" "1.73"
$ws.Range("E44").Value = "  +8.21%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "25.04"
$ws.Range("E45").Value = "  +7.97%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D46" "4.41"
$ws.Range("E46").Value = "  +0.08%  "

Set-TextValue "D47" "41.39"
$ws.Range("E47").Value = "  +0.31%  "

Set-TextValue "D48" "6.86"
$ws.Range("E48").Value = "  -1.73%  "

Set-TextValue "D49" "22.74"
$ws.Range("E49").Value = "  -0.33%  "

$ws.Range("D50").Value = "2.335.65"
$ws.Range("E50").Value = "  +5.45%  "

$ws.Range("E51").Value = "  +1.93%  "
